$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Mensalidade" values in column G, formatted with the built-in
# "Comma" cell style (thousands separator, 2 decimals).
$ws.Range("G2").Value = 150
$ws.Range("G3").Value = 75.5
$ws.Range("G2:G3").Style = "Comma"

# Widen column G to fit its new contents (best-fit for "150.00" / "75.50").
$ws.Columns.Item(7).ColumnWidth = 10

# Move the active selection, matching the author's final cursor position.
$ws.Range("I17").Select() | Out-Null
